# =====================================================================
# Reproduce the "Add files via upload" commit:
#  - new codeNames / workbookPr (not representable via this COM surface,
#    skipped - cosmetic VBA project metadata only)
#  - active tab moves from Sheet3 to Sheet2
#  - two new worksheets, Sheet4 (VLOOKUP demo) and Sheet5 (Text function
#    demo), appended after Sheet3
#  - two new cell styles (green fill, blue thin border) used on Sheet4
#  - a couple of volatile-formula rows added to Sheet1 / Sheet2
#  - TODAY() on Sheet3 recalculates to a newer cached date automatically
# =====================================================================

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# Sheet1: two new RANDBETWEEN rows under the existing formula block.
# ---------------------------------------------------------------------
$ws1.Range("L32").Formula = "=RANDBETWEEN(100,200)"
$ws1.Range("L33").Formula = "=RANDBETWEEN(100,200)"

# ---------------------------------------------------------------------
# Sheet2: one new RAND() row.
# ---------------------------------------------------------------------
$ws2.Range("I16").Formula = "=RAND()"

# ---------------------------------------------------------------------
# Sheet3: TODAY() cell is left alone - re-evaluating on save naturally
# refreshes its cached <v> the same way the source commit's did.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Add Sheet4 (appended after Sheet3) and Sheet5 (appended after Sheet4).
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Sheet4"

$ws5 = $wb.Worksheets.Add($null, $ws4)
$ws5.Name = "Sheet5"

# ---------------------------------------------------------------------
# Sheet4 content: small "name / salary" table plus a VLOOKUP demo.
# ---------------------------------------------------------------------
$ws4.Range("A1").Value = "name"
$ws4.Range("B1").Value = "salary"

$ws4.Range("A2").Value = "a"
$ws4.Range("B2").Value = 100
$ws4.Range("C2").Value = 130

$ws4.Range("A3").Value = "b"
$ws4.Range("B3").Value = 200
$ws4.Range("C3").Value = 260

$ws4.Range("A4").Value = "c"
$ws4.Range("B4").Value = 300
$ws4.Range("C4").Value = 390

$ws4.Range("A5").Value = "d"
$ws4.Range("B5").Value = 200
$ws4.Range("C5").Value = 260

$ws4.Range("D2").Formula = '=VLOOKUP("a",A2:B5,2,FALSE)'

# Green fill down column D (rows 1-10) and a blue thin box border down
# column E (rows 1-10) - two brand new cell styles.
$ws4.Range("D1:D10").Interior.Color = 65280      # RGB(0,255,0)
$ws4.Range("E1:E10").Borders.Color = 16711680    # RGB(0,0,255)
$ws4.Range("E1:E10").Borders.LineStyle = 1       # xlContinuous / thin

$ws4.Range("C2").Select()

# ---------------------------------------------------------------------
# Sheet5 content: "Text Function" demo sheet. Cell values are entered in
# the same order their text first appears in the shared-string table of
# the target workbook (53 salary already consumed above, then 54..74).
# ---------------------------------------------------------------------
$ws5.Range("A1").Value = "Text Function "
$ws5.Range("A3").Value = "proper"
$ws5.Range("B1").Value = "Text"
$ws5.Range("C1").Value = "Function "
$ws5.Range("B20").Value = "Deepak"
$ws5.Range("B3").Value = "deepak"
$ws5.Range("A4").Value = "upper"
$ws5.Range("A5").Value = "LOWER"
$ws5.Range("A6").Value = "left"
$ws5.Range("A7").Value = "right"
$ws5.Range("A8").Value = "mid"
$ws5.Range("A10").Value = "char"
$ws5.Range("B12").Value = "    deepak          "
$ws5.Range("B14").Value = "        deepak kumar    "
$ws5.Range("A17").Value = "trim&len"
$ws5.Range("A20").Value = "Concinate"
$ws5.Range("A27").Value = "Deepak kumar"
$ws5.Range("C27").Value = "First name"
$ws5.Range("D27").Value = "Last name"
$ws5.Range("D28").Value = "Kumar"
$ws5.Range("A28").Value = "raj kumar"

$ws5.Range("C3").Formula = "=PROPER(B3)"
$ws5.Range("C4").Formula = "=UPPER(B3)"
$ws5.Range("C5").Formula = "=LOWER(C4)"
$ws5.Range("C6").Formula = "=LEFT(B3,4)"
$ws5.Range("C7").Formula = "=RIGHT(B3,3)"
$ws5.Range("C8").Formula = "=MID(B3,2,4)"
$ws5.Range("C9").Formula = '=FIND("a",B3,1)'
$ws5.Range("C10").Formula = "=CHAR(85)"
$ws5.Range("C11").Formula = "=CHAR(102)"

$ws5.Range("C12").Formula = "=LEN(B12)"
$ws5.Range("C13").Formula = "=TRIM(B12)"
$ws5.Range("D13").Formula = "=LEN(C13)"
$ws5.Range("C14").Formula = "=LEN(B14)"
$ws5.Range("C15").Formula = "=TRIM(B14)"
$ws5.Range("D15").Formula = "=LEN(C15)"

$ws5.Range("C17").Formula = "=TRIM(LEN(B14))"
$ws5.Range("C18").Formula = "=LEN(TRIM(B14))"

$ws5.Range("C20").Formula = '="Hello"&B20'
$ws5.Range("C21").Formula = '="Hello"&"  "&B20'

$ws5.Range("C23").Formula = '=REPLACE(B3,2,1,"i")'
$ws5.Range("C24").Formula = "=REPT(B3,2)"

$ws5.Range("C28").Value = "Deepak"

# Column widths approximating the bestFit widths from the source file.
$ws5.Columns.Item(1).ColumnWidth = 12.7
$ws5.Columns.Item(2).ColumnWidth = 18.17
$ws5.Columns.Item(3).ColumnWidth = 12.86

$ws5.Range("B28").Select()

# ---------------------------------------------------------------------
# Make Sheet2 the final active / selected tab (activeTab moves 2 -> 1).
# ---------------------------------------------------------------------
$ws1.Range("E26").Select()
$ws2.Range("A6").Select()
